# Update the "Improved Time" figures on both sheets (Locations + Rewards).
# The chart series on each sheet read straight from these cells, so the
# embedded chart caches pick the new numbers up automatically on save.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Locations")
$ws2 = $wb.Worksheets.Item("Rewards")

# --- Locations sheet ---------------------------------------------------
$ws1.Range("D3").Value = 0.89
$ws1.Range("D8").Value = 203.755

# --- Rewards sheet -------------------------------------------------------
$ws2.Range("D3").Value = 1.677
$ws2.Range("D6").Value = 1026.41

# --- Restore the cursor/selection on each sheet, leaving "Rewards" active
# (matches the saved workbook state) ---------------------------------
$ws1.Range("D15").Select() | Out-Null
$ws2.Range("D6").Select() | Out-Null
